$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update valuation dates in column A (rows 3 and 5) from 2020-02-20 to 2021-02-20
$newDate = Get-Date -Year 2021 -Month 2 -Day 20 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("A3").Value = $newDate
$ws.Range("A5").Value = $newDate

# Move the active selection to A6
$ws.Range("A6").Select()
